# Apply updated TPM values to rows 2-5 (Ccl11-Ccr3 LR-pair sheet) and
# remove now-unused rows 6-9 (the dataset shrank from 8 data rows to 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Ccl11/Ccr3 -> Resolving-Mac -----------------------------
$rowVals = @('ECs', 'Ccl11', 'Ccr3', 'Resolving-Mac', 3, 1, 0.8417533333333332, 2.52526, 0.01079423211523897, 0.01079423211523897, 2, 0.6666666666666666, 0.1790523333333333, 0.537157, 1, 1, 0.1507178984244444, 1.35646108582, 0.01079423211523897, 0.01079423211523897)
for ($col = 1; $col -le $rowVals.Length; $col++) { $ws.Cells.Item(2, $col).Value = $rowVals[$col - 1] }

# --- Row 3: FAPs -> Ccl11/Ccr3 -> Resolving-Mac ----------------------------
$rowVals = @('FAPs', 'Ccl11', 'Ccr3', 'Resolving-Mac', 3, 1, 69.05064766666666, 207.151943, 0.885471656726338, 0.8854716567263378, 2, 0.6666666666666666, 0.1790523333333333, 0.537157, 1, 1, 12.36367958289456, 111.273116246051, 0.885471656726338, 0.8854716567263378)
for ($col = 1; $col -le $rowVals.Length; $col++) { $ws.Cells.Item(3, $col).Value = $rowVals[$col - 1] }

# --- Row 4: MuSCs -> Ccl11/Ccr3 -> Resolving-Mac ---------------------------
$rowVals = @('MuSCs', 'Ccl11', 'Ccr3', 'Resolving-Mac', 3, 1, 7.697976666666666, 23.09393, 0.09871507918910555, 0.09871507918910553, 2, 0.6666666666666666, 0.1790523333333333, 0.537157, 1, 1, 1.378340684112222, 12.40506615701, 0.09871507918910555, 0.09871507918910553)
for ($col = 1; $col -le $rowVals.Length; $col++) { $ws.Cells.Item(4, $col).Value = $rowVals[$col - 1] }

# --- Row 5: Resolving-Mac -> Ccl11/Ccr3 -> Resolving-Mac -------------------
$rowVals = @('Resolving-Mac', 'Ccl11', 'Ccr3', 'Resolving-Mac', 3, 1, 0.391393, 1.174179, 0.005019031969317685, 0.005019031969317684, 2, 0.6666666666666666, 0.1790523333333333, 0.537157, 1, 1, 0.07007982990033335, 0.630718469103, 0.005019031969317685, 0.005019031969317684)
for ($col = 1; $col -le $rowVals.Length; $col++) { $ws.Cells.Item(5, $col).Value = $rowVals[$col - 1] }

# --- Remove the old rows 6-9 (MuSCs/Resolving-Mac pairs no longer present) -
$ws.Range("A6:T9").EntireRow.Delete()
